# Natmi following Dr Hou advice
# Updates Ligand-expressing cells (E) and Receptor-expressing cells (K) from 1 to 3
# for every data row, and refreshes all of the downstream NATMI-derived metrics
# (average/total expression, derived specificity, and edge weights) that depend on
# those counts, matching the recomputed values from the updated pipeline run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 3.001642333333333
$ws.Range("N2").Value = 9.004927
$ws.Range("O2").Value = 0.05169795991651582
$ws.Range("P2").Value = 0.05169795991651582
$ws.Range("Q2").Value = 366.959825201635
$ws.Range("R2").Value = 3302.638426814715
$ws.Range("S2").Value = 0.01028085981454773
$ws.Range("T2").Value = 0.01028085981454773

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 11.05428166666667
$ws.Range("N3").Value = 33.162845
$ws.Range("O3").Value = 0.1903903753498087
$ws.Range("P3").Value = 0.1903903753498088
$ws.Range("Q3").Value = 1351.419262409225
$ws.Range("R3").Value = 12162.77336168303
$ws.Range("S3").Value = 0.03786177950099708
$ws.Range("T3").Value = 0.03786177950099708

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 16.49405866666666
$ws.Range("N4").Value = 49.482176
$ws.Range("O4").Value = 0.2840808761059341
$ws.Range("P4").Value = 0.2840808761059341
$ws.Range("Q4").Value = 2016.44840158688
$ws.Range("R4").Value = 18148.03561428192
$ws.Range("S4").Value = 0.05649344128772815
$ws.Range("T4").Value = 0.05649344128772816

# Row 5
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 27.51115433333333
$ws.Range("N5").Value = 82.533463
$ws.Range("O5").Value = 0.4738307886277414
$ws.Range("P5").Value = 0.4738307886277414
$ws.Range("Q5").Value = 3363.321563380315
$ws.Range("R5").Value = 30269.89407042284
$ws.Range("S5").Value = 0.09422785582960992
$ws.Range("T5").Value = 0.09422785582960992

# Row 6
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 3.001642333333333
$ws.Range("N6").Value = 9.004927
$ws.Range("O6").Value = 0.05169795991651582
$ws.Range("P6").Value = 0.05169795991651582
$ws.Range("Q6").Value = 397.8519676802439
$ws.Range("R6").Value = 3580.667709122195
$ws.Range("S6").Value = 0.01114634362062678
$ws.Range("T6").Value = 0.01114634362062678

# Row 7
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 11.05428166666667
$ws.Range("N7").Value = 33.162845
$ws.Range("O7").Value = 0.1903903753498087
$ws.Range("P7").Value = 0.1903903753498088
$ws.Range("Q7").Value = 1465.187128904536
$ws.Range("R7").Value = 13186.68416014083
$ws.Range("S7").Value = 0.04104913519094431
$ws.Range("T7").Value = 0.04104913519094432

# Row 8
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 16.49405866666666
$ws.Range("N8").Value = 49.482176
$ws.Range("O8").Value = 0.2840808761059341
$ws.Range("P8").Value = 0.2840808761059341
$ws.Range("Q8").Value = 2186.201074889351
$ws.Range("R8").Value = 19675.80967400416
$ws.Range("S8").Value = 0.06124928461855729
$ws.Range("T8").Value = 0.0612492846185573

# Row 9
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 27.51115433333333
$ws.Range("N9").Value = 82.533463
$ws.Range("O9").Value = 0.4738307886277414
$ws.Range("P9").Value = 0.4738307886277414
$ws.Range("Q9").Value = 3646.45939428655
$ws.Range("R9").Value = 32818.13454857895
$ws.Range("S9").Value = 0.1021603327598642
$ws.Range("T9").Value = 0.1021603327598642

# Row 10
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 3.001642333333333
$ws.Range("N10").Value = 9.004927
$ws.Range("O10").Value = 0.05169795991651582
$ws.Range("P10").Value = 0.05169795991651582
$ws.Range("Q10").Value = 960.793955525208
$ws.Range("R10").Value = 8647.145599726873
$ws.Range("S10").Value = 0.02691790024150976
$ws.Range("T10").Value = 0.02691790024150976

# Row 11
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3.0
$ws.Range("M11").Value = 11.05428166666667
$ws.Range("N11").Value = 33.162845
$ws.Range("O11").Value = 0.1903903753498087
$ws.Range("P11").Value = 0.1903903753498088
$ws.Range("Q11").Value = 3538.358614569488
$ws.Range("R11").Value = 31845.22753112539
$ws.Range("S11").Value = 0.09913174792362568
$ws.Range("T11").Value = 0.0991317479236257

# Row 12
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 16.49405866666666
$ws.Range("N12").Value = 49.482176
$ws.Range("O12").Value = 0.2840808761059341
$ws.Range("P12").Value = 0.2840808761059341
$ws.Range("Q12").Value = 5279.573683055344
$ws.Range("R12").Value = 47516.1631474981
$ws.Range("S12").Value = 0.1479141671332625
$ws.Range("T12").Value = 0.1479141671332626

# Row 13
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 27.51115433333333
$ws.Range("N13").Value = 82.533463
$ws.Range("O13").Value = 0.4738307886277414
$ws.Range("P13").Value = 0.4738307886277414
$ws.Range("Q13").Value = 8806.029452427922
$ws.Range("R13").Value = 79254.2650718513
$ws.Range("S13").Value = 0.2467124412691337
$ws.Range("T13").Value = 0.2467124412691338

# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3.0
$ws.Range("M14").Value = 3.001642333333333
$ws.Range("N14").Value = 9.004927
$ws.Range("O14").Value = 0.05169795991651582
$ws.Range("P14").Value = 0.05169795991651582
$ws.Range("Q14").Value = 119.6751596548178
$ws.Range("R14").Value = 1077.07643689336
$ws.Range("S14").Value = 0.003352856239831557
$ws.Range("T14").Value = 0.003352856239831558

# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3.0
$ws.Range("M15").Value = 11.05428166666667
$ws.Range("N15").Value = 33.162845
$ws.Range("O15").Value = 0.1903903753498087
$ws.Range("P15").Value = 0.1903903753498088
$ws.Range("Q15").Value = 440.7330309266222
$ws.Range("R15").Value = 3966.597278339601
$ws.Range("S15").Value = 0.01234771273424169
$ws.Range("T15").Value = 0.01234771273424169

# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3.0
$ws.Range("M16").Value = 16.49405866666666
$ws.Range("N16").Value = 49.482176
$ws.Range("O16").Value = 0.2840808761059341
$ws.Range("P16").Value = 0.2840808761059341
$ws.Range("Q16").Value = 657.6163596737421
$ws.Range("R16").Value = 5918.547237063679
$ws.Range("S16").Value = 0.01842398306638614
$ws.Range("T16").Value = 0.01842398306638614

# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3.0
$ws.Range("M17").Value = 27.51115433333333
$ws.Range("N17").Value = 82.533463
$ws.Range("O17").Value = 0.4738307886277414
$ws.Range("P17").Value = 0.4738307886277415
$ws.Range("Q17").Value = 1096.866788746871
$ws.Range("R17").Value = 9871.80109872184
$ws.Range("S17").Value = 0.030730158769133497
$ws.Range("T17").Value = 0.0307301587691335
